$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("w")

# Add new label and threshold value
$ws.Range("Q2").Value = "% unavail"
$ws.Range("R2").Value = 0.25

# Update the formulas for the availability grid to reference $R$2
# instead of the hardcoded 0.45 threshold. Set each shared-formula group
# separately (B2 alone, C2:O2, B3:O14) so Excel keeps the existing
# shared-formula boundaries instead of re-grouping everything together.
$ws.Range("B2").Formula = '=IF(RAND()> $R$2,1,0)'
$ws.Range("C2:O2").Formula = '=IF(RAND()> $R$2,1,0)'
$ws.Range("B3:O14").Formula = '=IF(RAND()> $R$2,1,0)'

# Update the selection shown in the sheet view
$ws.Range("A1:O14").Select()
